$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# A style-5 'donor' cell used to restore formatting on text-percentage cells
# after we nudge them with a leading apostrophe (prevents Excel's automatic
# text->percentage numeric coercion on assignment).
$styleDonor = $ws.Range("L5")

# --- Summary statistics (K/L columns) ---
$ws.Range("L6").Value = 207
$ws.Range("L7").Value = 9
$ws.Range("L9").Value = "'65.1%"
$styleDonor.Copy()
$ws.Range("L9").PasteSpecial(-4122)

# --- Recorded By column: swap order to 'System, dnasr281@gmail.com' ---
$ws.Range("G8").Value = "System, dnasr281@gmail.com"
$ws.Range("G9").Value = "System, dnasr281@gmail.com"
$ws.Range("G10").Value = "System, dnasr281@gmail.com"
$ws.Range("G12").Value = "System, dnasr281@gmail.com"
$ws.Range("G14").Value = "System, dnasr281@gmail.com"
$ws.Range("G15").Value = "System, dnasr281@gmail.com"
$ws.Range("G17").Value = "System, dnasr281@gmail.com"
$ws.Range("G34").Value = "System, dnasr281@gmail.com"
$ws.Range("G35").Value = "System, dnasr281@gmail.com"
$ws.Range("G36").Value = "System, dnasr281@gmail.com"
$ws.Range("G38").Value = "System, dnasr281@gmail.com"
$ws.Range("G40").Value = "System, dnasr281@gmail.com"
$ws.Range("G41").Value = "System, dnasr281@gmail.com"
$ws.Range("G43").Value = "System, dnasr281@gmail.com"
$ws.Range("G60").Value = "System, dnasr281@gmail.com"
$ws.Range("G61").Value = "System, dnasr281@gmail.com"
$ws.Range("G62").Value = "System, dnasr281@gmail.com"
$ws.Range("G64").Value = "System, dnasr281@gmail.com"
$ws.Range("G66").Value = "System, dnasr281@gmail.com"
$ws.Range("G67").Value = "System, dnasr281@gmail.com"
$ws.Range("G69").Value = "System, dnasr281@gmail.com"
$ws.Range("G86").Value = "System, dnasr281@gmail.com"
$ws.Range("G87").Value = "System, dnasr281@gmail.com"
$ws.Range("G88").Value = "System, dnasr281@gmail.com"
$ws.Range("G90").Value = "System, dnasr281@gmail.com"
$ws.Range("G92").Value = "System, dnasr281@gmail.com"
$ws.Range("G93").Value = "System, dnasr281@gmail.com"
$ws.Range("G95").Value = "System, dnasr281@gmail.com"
$ws.Range("G112").Value = "System, dnasr281@gmail.com"
$ws.Range("G113").Value = "System, dnasr281@gmail.com"
$ws.Range("G114").Value = "System, dnasr281@gmail.com"
$ws.Range("G116").Value = "System, dnasr281@gmail.com"
$ws.Range("G118").Value = "System, dnasr281@gmail.com"
$ws.Range("G119").Value = "System, dnasr281@gmail.com"
$ws.Range("G121").Value = "System, dnasr281@gmail.com"
$ws.Range("G138").Value = "System, dnasr281@gmail.com"
$ws.Range("G139").Value = "System, dnasr281@gmail.com"
$ws.Range("G140").Value = "System, dnasr281@gmail.com"
$ws.Range("G142").Value = "System, dnasr281@gmail.com"
$ws.Range("G144").Value = "System, dnasr281@gmail.com"
$ws.Range("G145").Value = "System, dnasr281@gmail.com"
$ws.Range("G147").Value = "System, dnasr281@gmail.com"

# --- Per-group class statistics (O/P/R/S columns), rows 15-20 ---
$ws.Range("O15").Value = 17
$ws.Range("P15").Value = 1
$ws.Range("R15").Value = "'65.4%"
$ws.Range("S15").Value = "'81.9%"
$styleDonor.Copy()
$ws.Range("R15:S15").PasteSpecial(-4122)
$ws.Range("O16").Value = 18
$ws.Range("P16").Value = 0
$ws.Range("R16").Value = "'69.2%"
$ws.Range("S16").Value = "'81.1%"
$styleDonor.Copy()
$ws.Range("R16:S16").PasteSpecial(-4122)
$ws.Range("O17").Value = 18
$ws.Range("P17").Value = 0
$ws.Range("R17").Value = "'69.2%"
$ws.Range("S17").Value = "'70.9%"
$styleDonor.Copy()
$ws.Range("R17:S17").PasteSpecial(-4122)
$ws.Range("O18").Value = 18
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = "'69.2%"
$ws.Range("S18").Value = "'77.0%"
$styleDonor.Copy()
$ws.Range("R18:S18").PasteSpecial(-4122)
$ws.Range("O19").Value = 18
$ws.Range("P19").Value = 0
$ws.Range("R19").Value = "'69.2%"
$ws.Range("S19").Value = "'76.5%"
$styleDonor.Copy()
$ws.Range("R19:S19").PasteSpecial(-4122)
$ws.Range("O20").Value = 17
$ws.Range("P20").Value = 1
$ws.Range("R20").Value = "'65.4%"
$ws.Range("S20").Value = "'80.1%"
$styleDonor.Copy()
$ws.Range("R20:S20").PasteSpecial(-4122)

# --- Rows that flip from 'Not Recorded' to 'Recorded' (style + values) ---
$ws.Range("A18:I18").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)
$ws.Range("G19").Value = "dnasr281@gmail.com"
$ws.Range("H19").Value = "19/26"
$ws.Range("I19").Value = "Recorded"
$ws.Range("A44:I44").Copy()
$ws.Range("A45:I45").PasteSpecial(-4122)
$ws.Range("G45").Value = "dnasr281@gmail.com"
$ws.Range("H45").Value = "20/27"
$ws.Range("I45").Value = "Recorded"
$ws.Range("A70:I70").Copy()
$ws.Range("A71:I71").PasteSpecial(-4122)
$ws.Range("G71").Value = "dnasr281@gmail.com"
$ws.Range("H71").Value = "17/26"
$ws.Range("I71").Value = "Recorded"
$ws.Range("A96:I96").Copy()
$ws.Range("A97:I97").PasteSpecial(-4122)
$ws.Range("G97").Value = "dnasr281@gmail.com"
$ws.Range("H97").Value = "24/27"
$ws.Range("I97").Value = "Recorded"
$ws.Range("A122:I122").Copy()
$ws.Range("A123:I123").PasteSpecial(-4122)
$ws.Range("G123").Value = "dnasr281@gmail.com"
$ws.Range("H123").Value = "25/30"
$ws.Range("I123").Value = "Recorded"
$ws.Range("A148:I148").Copy()
$ws.Range("A149:I149").PasteSpecial(-4122)
$ws.Range("G149").Value = "dnasr281@gmail.com"
$ws.Range("H149").Value = "17/23"
$ws.Range("I149").Value = "Recorded"

$excel.CutCopyMode = $false
